$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "K" = 7.94
    "L" = 0.2391566265060241
    "M" = 2.31
    "N" = 0.04529411764705882
    "O" = 0.2909319899244333
    "P" = 2.31
    "Q" = 0.04529411764705882
    "R" = 0.2909319899244333
    "T" = 0
    "U" = 48.7
    "V" = 0.9549019607843138
    "W" = 0.2193370165745856
    "X" = 0.09563064629785738
    "Y" = 0.1237063702767282
    "Z" = -3.223300970873792
    "AB" = 0.07579453389138853
    "AC" = -0.07579453389138853
    "AD" = 42.7
    "AF" = 42.7
    "AG" = -6
    "AH" = 0.455709711846318
    "AI" = 0.526510480887793
    "AJ" = -0.1333333333333333
    "AK" = -0.1851851851851852
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range($col + "2").Value = $val
    $ws.Range($col + "3").Value = $val
}
